# updated the dose_response script
# - fluorescence sheet: header row had the "+ SENSOR" / "- SENSOR" columns
#   swapped relative to the od600 sheet; fix the header labels so both
#   sheets agree (B:D -> "- SENSOR_x", E:G -> "+ SENSOR_x").
# - fluorescence sheet: correct a data entry (C4) from 185 to 197.
# - metadata sheet: new chart title, new (shorter) color list.

$wb = $excel.ActiveWorkbook

$fluor = $wb.Worksheets.Item("fluorescence")
$fluor.Range("B1").Value = " - SENSOR_1"
$fluor.Range("C1").Value = " - SENSOR_2"
$fluor.Range("D1").Value = " - SENSOR_3"
$fluor.Range("E1").Value = " + SENSOR_1"
$fluor.Range("F1").Value = " + SENSOR_2"
$fluor.Range("G1").Value = " + SENSOR_3"
$fluor.Range("C4").Value = 197.0

$meta = $wb.Worksheets.Item("metadata")
$meta.Range("A2").Value = "Promoter activity with and without cognate biosensor"
$meta.Range("D2").Value = "#6e6e6e"
$meta.Range("D3").Value = "#1fde0d"
$meta.Range("D4:D7").Clear()
